# Applies the "added expert user data" change:
#  - Copies Expert 1 / Expert 2 summary tables into the "All Data" sheet
#    (columns I:N, rows 1-10)
#  - Fills in the missing measurement values on the "Expert 1" and
#    "Expert 2" sheets (columns C:F)
#  - Leaves "All Data" as the active sheet / selected cell, and updates
#    the selection on the two expert sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Expert 1" sheet (already has headers + trial numbers) - add the
#    measured values for trials 1-3.
# ---------------------------------------------------------------------
$wsExpert1 = $wb.Worksheets.Item("Expert 1")

$wsExpert1.Range("C3").Value = 9.2
$wsExpert1.Range("D3").Value = 19.8
$wsExpert1.Range("E3").Value = 49.7
$wsExpert1.Range("F3").Value = 194.1

$wsExpert1.Range("C4").Value = 9.8
$wsExpert1.Range("D4").Value = 20.5
$wsExpert1.Range("E4").Value = 49.1
$wsExpert1.Range("F4").Value = 198.8

$wsExpert1.Range("C5").Value = 10
$wsExpert1.Range("D5").Value = 20.3
$wsExpert1.Range("E5").Value = 50.6
$wsExpert1.Range("F5").Value = 200.4

# ---------------------------------------------------------------------
# 2. "Expert 2" sheet was completely empty - add the header row, the
#    "Trial" label, the trial numbers and the measured values.
# ---------------------------------------------------------------------
$wsExpert2 = $wb.Worksheets.Item("Expert 2")

$wsExpert2.Range("B1").Value = "Measurement "
$wsExpert2.Range("C1").Value = 10
$wsExpert2.Range("D1").Value = 20
$wsExpert2.Range("E1").Value = 50
$wsExpert2.Range("F1").Value = 200
$wsExpert2.Range("J1").Value = "weight in mg"

$wsExpert2.Range("A2").Value = "Trial"

$wsExpert2.Range("A3").Value = 1
$wsExpert2.Range("C3").Value = 10.3
$wsExpert2.Range("D3").Value = 20.3
$wsExpert2.Range("E3").Value = 49.3
$wsExpert2.Range("F3").Value = 192.6

$wsExpert2.Range("A4").Value = 2
$wsExpert2.Range("C4").Value = 10.3
$wsExpert2.Range("D4").Value = 20.6
$wsExpert2.Range("E4").Value = 49.1
$wsExpert2.Range("F4").Value = 192.5

$wsExpert2.Range("A5").Value = 3
$wsExpert2.Range("C5").Value = 9.1
$wsExpert2.Range("D5").Value = 18.9
$wsExpert2.Range("E5").Value = 48.7
$wsExpert2.Range("F5").Value = 191.4

# ---------------------------------------------------------------------
# 3. "All Data" sheet - paste in the two expert tables side-by-side in
#    columns I:N (Expert 2's table on rows 1-5, Expert 1's on rows 6-10).
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Data")

# Expert 2 block -> rows 1-5
$wsAll.Range("J1").Value = "Measurement "
$wsAll.Range("K1").Value = 10
$wsAll.Range("L1").Value = 20
$wsAll.Range("M1").Value = 50
$wsAll.Range("N1").Value = 200

$wsAll.Range("I2").Value = "Trial"

$wsAll.Range("I3").Value = 1
$wsAll.Range("K3").Value = 10.3
$wsAll.Range("L3").Value = 20.3
$wsAll.Range("M3").Value = 49.3
$wsAll.Range("N3").Value = 192.6

$wsAll.Range("I4").Value = 2
$wsAll.Range("K4").Value = 10.3
$wsAll.Range("L4").Value = 20.6
$wsAll.Range("M4").Value = 49.1
$wsAll.Range("N4").Value = 192.5

$wsAll.Range("I5").Value = 3
$wsAll.Range("K5").Value = 9.1
$wsAll.Range("L5").Value = 18.9
$wsAll.Range("M5").Value = 48.7
$wsAll.Range("N5").Value = 191.4

# Expert 1 block -> rows 6-10
$wsAll.Range("J6").Value = "Measurement "
$wsAll.Range("K6").Value = 10
$wsAll.Range("L6").Value = 20
$wsAll.Range("M6").Value = 50
$wsAll.Range("N6").Value = 200

$wsAll.Range("I7").Value = "Trial"

$wsAll.Range("I8").Value = 1
$wsAll.Range("K8").Value = 9.2
$wsAll.Range("L8").Value = 19.8
$wsAll.Range("M8").Value = 49.7
$wsAll.Range("N8").Value = 194.1

$wsAll.Range("I9").Value = 2
$wsAll.Range("K9").Value = 9.8
$wsAll.Range("L9").Value = 20.5
$wsAll.Range("M9").Value = 49.1
$wsAll.Range("N9").Value = 198.8

$wsAll.Range("I10").Value = 3
$wsAll.Range("K10").Value = 10
$wsAll.Range("L10").Value = 20.3
$wsAll.Range("M10").Value = 50.6
$wsAll.Range("N10").Value = 200.4

# Widen column J (10) on "All Data" to fit the new "Measurement " header
$wsAll.Columns.Item(10).ColumnWidth = 10.833333333333334

# ---------------------------------------------------------------------
# 4. Selections / active sheet - "All Data" becomes the active tab, the
#    two expert sheets keep a selection over their newly-filled table.
# ---------------------------------------------------------------------
$wsExpert1.Range("A1:F5").Select()
$wsExpert1.Range("F5").Activate()

$wsExpert2.Range("A1:F6").Select()
$wsExpert2.Range("F6").Activate()

$wsAll.Activate()
$wsAll.Range("D36").Select()
